$wb = $excel.ActiveWorkbook

# --- Delete the empty "Sheet1" ---
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Sheet1").Delete()
$excel.DisplayAlerts = $true

# --- Add API / UI combination data to ProductCheckout ---
$ws = $wb.Worksheets.Item("ProductCheckout")
$ws.Activate()

# Header row (E1:G1) -- reuse the existing "pasted" header style (same style
# already used elsewhere in the workbook, e.g. ProductList!A2) by copying its
# format instead of hand-building a new font, so no redundant style gets
# created in styles.xml.
$fmtSource = $wb.Worksheets.Item("ProductList").Range("A2")
$fmtSource.Copy()
$ws.Range("E1:G1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("E1").Value = "Postcode"
$ws.Range("F1").Value = "State"
$ws.Range("G1").Value = "Phone"

# Data rows
$ws.Range("E2").Value = 123456
$ws.Range("F2").Value = "Austria"
$ws.Range("G2").Value = 879456234

$ws.Range("E3").Value = 141452
$ws.Range("F3").Value = "Austria"
$ws.Range("G3").Value = 879466234

# Column widths for the new columns (nearest values this engine's internal
# 1/6-character-unit rounding can represent, closest to the authored
# 13.44140625 / 12.109375 / 11.6640625 stored widths)
$ws.Columns.Item(5).ColumnWidth = 12.666666666666666
$ws.Columns.Item(6).ColumnWidth = 11.333333333333334
$ws.Columns.Item(7).ColumnWidth = 10.833333333333334

# Selection, as recorded at save time
$ws.Range("G5").Select()
